$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 37-41 (re-ordered match results for matchday) ---
# Row 37: Chippa Utd. vs Royal AM
$ws.Range("F37").Value = "Chippa Utd."
$ws.Range("G37").Value = 2
$ws.Range("H37").Value = "Royal AM"
$ws.Range("I37").Value = 3
$ws.Range("J37").Value = 2.24
$ws.Range("K37").Value = "28/08/2023 13:24"
$ws.Range("L37").Value = 2.16
$ws.Range("M37").Value = "30/08/2023 19:27"
$ws.Range("N37").Value = 3
$ws.Range("O37").Value = "28/08/2023 13:24"
$ws.Range("P37").Value = 3.01
$ws.Range("Q37").Value = "30/08/2023 19:27"
$ws.Range("R37").Value = 3.42
$ws.Range("S37").Value = "28/08/2023 13:24"
$ws.Range("T37").Value = 4
$ws.Range("U37").Value = "30/08/2023 19:27"
$ws.Range("V37").Value = "https://www.betexplorer.com/football/south-africa/premier-league/chippa-utd-royal-am/SOPo3lUs/"

# Row 38: AmaZulu vs TS Galaxy
$ws.Range("F38").Value = "AmaZulu"
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = "TS Galaxy"
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 2.58
$ws.Range("K38").Value = "28/08/2023 13:23"
$ws.Range("L38").Value = 3.23
$ws.Range("M38").Value = "30/08/2023 19:19"
$ws.Range("N38").Value = 2.77
$ws.Range("O38").Value = "28/08/2023 13:23"
$ws.Range("P38").Value = 2.75
$ws.Range("Q38").Value = "30/08/2023 19:19"
$ws.Range("R38").Value = 3.1
$ws.Range("S38").Value = "28/08/2023 13:23"
$ws.Range("T38").Value = 2.7
$ws.Range("U38").Value = "30/08/2023 19:19"
$ws.Range("V38").Value = "https://www.betexplorer.com/football/south-africa/premier-league/amazulu-ts-galaxy/0xOg1Sbf/"

# Row 39: Richards Bay vs Sekhukhune
$ws.Range("F39").Value = "Richards Bay"
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = "Sekhukhune"
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 2.96
$ws.Range("K39").Value = "28/08/2023 13:24"
$ws.Range("L39").Value = 3.28
$ws.Range("M39").Value = "30/08/2023 19:29"
$ws.Range("N39").Value = 2.79
$ws.Range("O39").Value = "28/08/2023 13:24"
$ws.Range("P39").Value = 2.6
$ws.Range("Q39").Value = "30/08/2023 19:29"
$ws.Range("R39").Value = 2.67
$ws.Range("S39").Value = "28/08/2023 13:24"
$ws.Range("T39").Value = 2.82
$ws.Range("U39").Value = "30/08/2023 19:29"
$ws.Range("V39").Value = "https://www.betexplorer.com/football/south-africa/premier-league/richards-bay-sekhukhune/djxUp7ED/"

# Row 40: Swallows vs Cape Town Spurs
$ws.Range("F40").Value = "Swallows"
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = "Cape Town Spurs"
$ws.Range("I40").Value = 1
$ws.Range("J40").Value = 1.76
$ws.Range("K40").Value = "28/08/2023 13:24"
$ws.Range("L40").Value = 1.83
$ws.Range("M40").Value = "30/08/2023 19:26"
$ws.Range("N40").Value = 3.23
$ws.Range("O40").Value = "28/08/2023 13:24"
$ws.Range("P40").Value = 3.23
$ws.Range("Q40").Value = "30/08/2023 19:26"
$ws.Range("R40").Value = 5.01
$ws.Range("S40").Value = "28/08/2023 13:24"
$ws.Range("T40").Value = 5.17
$ws.Range("U40").Value = "30/08/2023 19:26"
$ws.Range("V40").Value = "https://www.betexplorer.com/football/south-africa/premier-league/swallows-fc-cape-town-spurs/MoPk28ql/"

# Row 41: Stellenbosch vs Kaizer Chiefs
$ws.Range("F41").Value = "Stellenbosch"
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = "Kaizer Chiefs"
$ws.Range("I41").Value = 2
$ws.Range("J41").Value = 2.81
$ws.Range("K41").Value = "28/08/2023 13:24"
$ws.Range("L41").Value = 2.68
$ws.Range("M41").Value = "30/08/2023 19:26"
$ws.Range("N41").Value = 2.8
$ws.Range("O41").Value = "28/08/2023 13:24"
$ws.Range("P41").Value = 2.73
$ws.Range("Q41").Value = "30/08/2023 19:26"
$ws.Range("R41").Value = 2.81
$ws.Range("S41").Value = "28/08/2023 13:24"
$ws.Range("T41").Value = 3.29
$ws.Range("U41").Value = "30/08/2023 19:26"
$ws.Range("V41").Value = "https://www.betexplorer.com/football/south-africa/premier-league/stellenbosch-fc-kaizer-chiefs/6LJK7ULQ/"

# --- Update rows 63-65 (re-ordered match results for matchday) ---
# Row 63: Cape Town Spurs vs Supersport Utd
$ws.Range("F63").Value = "Cape Town Spurs"
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = "Supersport Utd"
$ws.Range("I63").Value = 2
$ws.Range("J63").Value = 4.52
$ws.Range("K63").Value = "30/09/2023 14:13"
$ws.Range("L63").Value = 5.17
$ws.Range("M63").Value = "04/10/2023 19:21"
$ws.Range("N63").Value = 3.07
$ws.Range("O63").Value = "30/09/2023 14:13"
$ws.Range("P63").Value = 3.52
$ws.Range("Q63").Value = "04/10/2023 19:21"
$ws.Range("R63").Value = 1.96
$ws.Range("S63").Value = "30/09/2023 14:13"
$ws.Range("T63").Value = 1.75
$ws.Range("U63").Value = "04/10/2023 19:21"
$ws.Range("V63").Value = "https://www.betexplorer.com/football/south-africa/premier-league/cape-town-spurs-supersport-utd/fHreMlhH/"

# Row 64: Richards Bay vs Polokwane
$ws.Range("F64").Value = "Richards Bay"
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = "Polokwane"
$ws.Range("I64").Value = 1
$ws.Range("J64").Value = 2.58
$ws.Range("K64").Value = "03/10/2023 01:12"
$ws.Range("L64").Value = 2.42
$ws.Range("M64").Value = "04/10/2023 19:21"
$ws.Range("N64").Value = 2.85
$ws.Range("O64").Value = "03/10/2023 01:12"
$ws.Range("P64").Value = 2.83
$ws.Range("Q64").Value = "04/10/2023 19:21"
$ws.Range("R64").Value = 3.17
$ws.Range("S64").Value = "03/10/2023 01:12"
$ws.Range("T64").Value = 3.61
$ws.Range("U64").Value = "04/10/2023 19:21"
$ws.Range("V64").Value = "https://www.betexplorer.com/football/south-africa/premier-league/richards-bay-polokwane-city/4SkrPj8b/"

# Row 65: Stellenbosch vs TS Galaxy
$ws.Range("F65").Value = "Stellenbosch"
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = "TS Galaxy"
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2.27
$ws.Range("K65").Value = "03/10/2023 01:12"
$ws.Range("L65").Value = 2.3
$ws.Range("M65").Value = "04/10/2023 19:22"
$ws.Range("N65").Value = 2.99
$ws.Range("O65").Value = "03/10/2023 01:12"
$ws.Range("P65").Value = 2.93
$ws.Range("Q65").Value = "04/10/2023 19:22"
$ws.Range("R65").Value = 3.58
$ws.Range("S65").Value = "03/10/2023 01:12"
$ws.Range("T65").Value = 3.74
$ws.Range("U65").Value = "04/10/2023 19:22"
$ws.Range("V65").Value = "https://www.betexplorer.com/football/south-africa/premier-league/stellenbosch-fc-ts-galaxy/niumOAN4/"

# --- Append new rows 75-78, copying formatting from row 74 first ---
$ws.Range("A74:V74").Copy()
$ws.Range("A75:V78").PasteSpecial(-4122)

# Row 75: Cape Town City vs Chippa Utd.
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = "south-africa"
$ws.Range("C75").Value = "premier-league"
$ws.Range("D75").Value = "2023-2024"
$ws.Range("E75").Value = 45237.77083333334
$ws.Range("F75").Value = "Cape Town City"
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = "Chippa Utd."
$ws.Range("I75").Value = 1
$ws.Range("J75").Value = 1.71
$ws.Range("K75").Value = "06/11/2023 15:19"
$ws.Range("L75").Value = 1.83
$ws.Range("M75").Value = "07/11/2023 18:28"
$ws.Range("N75").Value = 3.53
$ws.Range("O75").Value = "06/11/2023 15:19"
$ws.Range("P75").Value = 3.19
$ws.Range("Q75").Value = "07/11/2023 18:28"
$ws.Range("R75").Value = 5.54
$ws.Range("S75").Value = "06/11/2023 15:19"
$ws.Range("T75").Value = 5.37
$ws.Range("U75").Value = "07/11/2023 18:28"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/south-africa/premier-league/cape-town-city-chippa-utd/S8L2Iuiq/"

# Row 76: Orlando Pirates vs Sekhukhune
$ws.Range("A76").Value = 75
$ws.Range("B76").Value = "south-africa"
$ws.Range("C76").Value = "premier-league"
$ws.Range("D76").Value = "2023-2024"
$ws.Range("E76").Value = 45237.77083333334
$ws.Range("F76").Value = "Orlando Pirates"
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = "Sekhukhune"
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 1.62
$ws.Range("K76").Value = "06/11/2023 13:00"
$ws.Range("L76").Value = 1.66
$ws.Range("M76").Value = "07/11/2023 18:22"
$ws.Range("N76").Value = 3.8
$ws.Range("O76").Value = "06/11/2023 13:00"
$ws.Range("P76").Value = 3.55
$ws.Range("Q76").Value = "07/11/2023 18:22"
$ws.Range("R76").Value = 5.04
$ws.Range("S76").Value = "06/11/2023 13:00"
$ws.Range("T76").Value = 6.06
$ws.Range("U76").Value = "07/11/2023 18:22"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/south-africa/premier-league/orlando-pirates-sekhukhune/6aK6Ha6k/"

# Row 77: Royal AM vs Golden Arrows
$ws.Range("A77").Value = 76
$ws.Range("B77").Value = "south-africa"
$ws.Range("C77").Value = "premier-league"
$ws.Range("D77").Value = "2023-2024"
$ws.Range("E77").Value = 45237.77083333334
$ws.Range("F77").Value = "Royal AM"
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = "Golden Arrows"
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3.64
$ws.Range("K77").Value = "06/11/2023 15:19"
$ws.Range("L77").Value = 3.25
$ws.Range("M77").Value = "07/11/2023 18:26"
$ws.Range("N77").Value = 3.19
$ws.Range("O77").Value = "06/11/2023 15:19"
$ws.Range("P77").Value = 3.07
$ws.Range("Q77").Value = "07/11/2023 18:26"
$ws.Range("R77").Value = 2.19
$ws.Range("S77").Value = "06/11/2023 15:19"
$ws.Range("T77").Value = 2.43
$ws.Range("U77").Value = "07/11/2023 18:26"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/south-africa/premier-league/royal-am-golden-arrows/EPJAGJLe/"

# Row 78: Stellenbosch vs Richards Bay
$ws.Range("A78").Value = 77
$ws.Range("B78").Value = "south-africa"
$ws.Range("C78").Value = "premier-league"
$ws.Range("D78").Value = "2023-2024"
$ws.Range("E78").Value = 45237.77083333334
$ws.Range("F78").Value = "Stellenbosch"
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = "Richards Bay"
$ws.Range("I78").Value = 1
$ws.Range("J78").Value = 2.03
$ws.Range("K78").Value = "06/11/2023 15:19"
$ws.Range("L78").Value = 1.99
$ws.Range("M78").Value = "07/11/2023 18:22"
$ws.Range("N78").Value = 3.24
$ws.Range("O78").Value = "06/11/2023 15:19"
$ws.Range("P78").Value = 3.17
$ws.Range("Q78").Value = "07/11/2023 18:22"
$ws.Range("R78").Value = 4.1
$ws.Range("S78").Value = "06/11/2023 15:19"
$ws.Range("T78").Value = 4.37
$ws.Range("U78").Value = "07/11/2023 18:21"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/south-africa/premier-league/stellenbosch-fc-richards-bay/hGIEFwy2/"
